$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 data updates
$ws.Range("C8").Value = "'"
$ws.Range("F8").Value = "'0.300"
$ws.Range("I8").Value = "RS800"
$ws.Range("L8").Value = 0.001
$ws.Range("M8").Value = "'0.008"
$ws.Range("G8").Value = "'0.612"

# Row 4: set test case id reference
$ws.Range("B4").Value = "NGC-488/T396 OR TC-149"

# Update sheet view selection
$null = $ws.Range("B7").Select()
